$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

Set-TextValue $ws "D2" "22.395.38"
Set-TextValue $ws "E2" "  -4.56%  "
Set-TextValue $ws "D3" "1.570.31"
Set-TextValue $ws "E3" "  -4.70%  "
Set-TextValue $ws "D4" "1.002"
Set-TextValue $ws "E4" "  +0.11%  "
Set-TextValue $ws "E5" "  -0.04%  "
Set-TextValue $ws "D6" "291.70"
Set-TextValue $ws "E6" "  -2.30%  "
Set-TextValue $ws "D7" "0.3663"
Set-TextValue $ws "E7" "  -3.14%  "
Set-TextValue $ws "D8" "49.36"
Set-TextValue $ws "E8" "  -1.08%  "
Set-TextValue $ws "E9" "  -4.78%  "
Set-TextValue $ws "E10" "  -3.79%  "
Set-TextValue $ws "D11" "0.07585"
Set-TextValue $ws "E11" "  -6.27%  "
Set-TextValue $ws "D12" "0.9982"
Set-TextValue $ws "E12" "  -0.25%  "
Set-TextValue $ws "D13" "21.15"
Set-TextValue $ws "E13" "  -4.20%  "
Set-TextValue $ws "D14" "6.071"
Set-TextValue $ws "E14" "  -4.99%  "
Set-TextValue $ws "D15" "6.883"
Set-TextValue $ws "E15" "  -6.35%  "
Set-TextValue $ws "D16" "0.00001140"
Set-TextValue $ws "E16" "  -4.78%  "
Set-TextValue $ws "D17" "1.571.68"
Set-TextValue $ws "E17" "  -5.17%  "
Set-TextValue $ws "D18" "89.02"
Set-TextValue $ws "D19" "0.06720"
Set-TextValue $ws "D20" "0.9993"
Set-TextValue $ws "E20" "  -0.07%  "
Set-TextValue $ws "D21" "6.270"
Set-TextValue $ws "E21" "  -7.25%  "
Set-TextValue $ws "D22" "16.46"
Set-TextValue $ws "E22" "  -5.17%  "
Set-TextValue $ws "D23" "0.5253"
Set-TextValue $ws "E23" "  -8.77%  "
Set-TextValue $ws "E24" "  -3.18%  "
Set-TextValue $ws "D25" "22.467.95"
Set-TextValue $ws "E25" "  -4.24%  "
Set-TextValue $ws "E28" "  -4.72%  "
Set-TextValue $ws "D29" "144.61"
Set-TextValue $ws "E29" "  -5.67%  "
Set-TextValue $ws "D30" "4.998"
Set-TextValue $ws "E30" "  -3.72%  "
Set-TextValue $ws "D31" "125.15"
Set-TextValue $ws "E31" "  -5.84%  "
Set-TextValue $ws "D32" "1.748.13"
Set-TextValue $ws "E32" "  -4.80%  "
Set-TextValue $ws "D33" "1.048"
Set-TextValue $ws "E33" "  +5.60%  "
Set-TextValue $ws "D34" "6.289"
Set-TextValue $ws "E34" "  -9.30%  "
Set-TextValue $ws "D35" "1.980"
Set-TextValue $ws "E35" "  -7.36%  "
Set-TextValue $ws "D36" "10.44"
Set-TextValue $ws "E36" "  -8.86%  "
Set-TextValue $ws "D37" "0.02569"
Set-TextValue $ws "E37" "  -5.50%  "
Set-TextValue $ws "D38" "0.08455"
Set-TextValue $ws "E38" "  -3.23%  "
Set-TextValue $ws "D39" "0.2304"
Set-TextValue $ws "E39" "  -5.23%  "
Set-TextValue $ws "D40" "0.06532"
Set-TextValue $ws "E40" "  -3.80%  "
Set-TextValue $ws "D41" "5.524"
Set-TextValue $ws "E41" "  -6.89%  "
Set-TextValue $ws "D42" "11.86"
Set-TextValue $ws "E42" "  -9.09%  "
Set-TextValue $ws "D43" "1.253"
Set-TextValue $ws "E43" "  -3.41%  "
Set-TextValue $ws "D44" "0.6409"
Set-TextValue $ws "E44" "  -7.10%  "
Set-TextValue $ws "D45" "14.60"
Set-TextValue $ws "E45" "  -6.74%  "
Set-TextValue $ws "D46" "0.9992"
Set-TextValue $ws "E46" "  -0.10%  "
Set-TextValue $ws "D47" "0.6029"
Set-TextValue $ws "E47" "  -5.24%  "
Set-TextValue $ws "D48" "3.784"
Set-TextValue $ws "E48" "  -3.14%  "
Set-TextValue $ws "D49" "2.139"
Set-TextValue $ws "E49" "  -5.08%  "
Set-TextValue $ws "D50" "123.29"
Set-TextValue $ws "E50" "  -3.30%  "
Set-TextValue $ws "D51" "1.211"
Set-TextValue $ws "E51" "  +2.41%  "

# Row 26/27 swap: Toncoin <-> LidoDAOToken with updated values
Set-TextValue $ws "B26" "LidoDAOToken"
Set-TextValue $ws "C26" "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue $ws "D26" "3.025"
Set-TextValue $ws "E26" "  +4.45%  "
Set-TextValue $ws "B27" "Toncoin"
Set-TextValue $ws "C27" "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue $ws "D27" "2.384"
Set-TextValue $ws "E27" "  -4.22%  "
